$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Counter cell update
$ws.Range("F1").Value = 9

# Row 4
$ws.Range("A4").Value = "camilo"
$ws.Range("B4").Value = "37V9#A1>P0%j4-i3/i4#X3[c4)g6%F2|b9)g9+i2``X2`$T7)E0(s5-U6(K3]Y6;,87V9#90A1>139P0%92j4-73i3/158i4#93X3[163c4)159g6%86F2|101b9)123g9+91i2``79X2`$154T7)88E0(74s5-108U6(89K3]162Y6;,124f7`$122y8~134z1@130V6*133u3%136G4},"
$ws.Range("C4").Value = "51m0>z2:n6]j2{L0?p4}l4@j8/B7}u2?o1;N7[S1:F3`$t0&Q9{u9;X7=F4}l0*(108m0>101z2:156n6]107j2{94L0?133p4}100l4@92j8/153B7}99u2?142o1;136N7[101S1:84F3`$125t0&107Q9{97u9;170X7=106F4}109l0*(191i8+189R1)201f0[197U2!200J8.203w7:("

# Row 5
$ws.Range("A5").Value = "camilo"
$ws.Range("B5").Value = "35b6!T3:o1)J6^T3]H3`$n8{o8<;87b6!89T3:143o1)83J6^76T3]120H3`$86n8{78o8<;166l0)166U3+;"
$ws.Range("C5").Value = "10M1~N2[U7}F5|K1;q0=H1+h4.y4=j3.y3.*65M1~65N2[79U7}62F5|55K1;113q0=62H1+69h4.92y4=61j3.105y3.*197E4-197g4;197R3_*"

# Row 6
$ws.Range("A6").Value = "ññ"
$ws.Range("B6").Value = "76k1@J0``q1}Z0/Z9&m2@g1/t7>{132k1@130J0``146q1}127Z0/114Z9&153m2@127g1/171t7>{327F3&327M3_{"
$ws.Range("C6").Value = "11m1}L6|R0[w0*O2?)60m1}65L6|132R0[68w0*57O2?)257y9.)"

# Row 7
$ws.Range("A7").Value = "ñ"
$ws.Range("B7").Value = "62d9`$M5<R2{V4(g5!&116d9`$110M5<128R2{118V4(125g5!&301B8?&"
$ws.Range("C7").Value = "60d7%O5/T8}L9*L8>:115d7%111O5/162T8}108L9*106L8>:314f0.:"

# Row 8
$ws.Range("A8").Value = "l"
$ws.Range("B8").Value = "46m1+J3)b8!J1,w5,-101m1+101J3)158b8!97J1,108w5,-185p3>-"
$ws.Range("C8").Value = "95B2*H2^P7.R1&P6|?146B2*151H2^177P7.149R1&156P6|?146R6=?"

# Clear the style that the column-level default (applied to columns B and C)
# would otherwise stamp onto these freshly-written cells, matching the
# target which leaves them with no explicit cell style.
$ws.Range("A4:C8").Style = "Normal"
